# Auto-generated PowerShell Excel COM-interop script
# Updates run_time, max_er, and iteration columns for rows 2-11 (run_number 0-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 3).Value = 0.821000337600708
$ws.Cells.Item(2, 5).Value = 1625.507514015328
$ws.Cells.Item(2, 6).Value = 0.1763624076937847
$ws.Cells.Item(2, 7).Value = 0.1313833334672443
$ws.Cells.Item(2, 8).Value = 0.09292159289432228
$ws.Cells.Item(2, 9).Value = 0.07328782157931228
$ws.Cells.Item(2, 10).Value = 0.06496955705674243
$ws.Cells.Item(2, 11).Value = 0.05707105522623642
$ws.Cells.Item(2, 12).Value = 0.05040991018545042
$ws.Cells.Item(2, 13).Value = 0.04714180182084696
$ws.Cells.Item(2, 14).Value = 0.04279921585802832
$ws.Cells.Item(2, 15).Value = 0.04064702430601397
$ws.Cells.Item(2, 16).Value = 0.0379955077036014
$ws.Cells.Item(2, 17).Value = 0.03647308144372229
$ws.Cells.Item(2, 18).Value = 0.0350255879848597
$ws.Cells.Item(2, 19).Value = 0.03425499183729432
$ws.Cells.Item(2, 20).Value = 0.03306593189516139
$ws.Cells.Item(2, 21).Value = 0.03261324237549594
$ws.Cells.Item(2, 22).Value = 0.03220966239773476
$ws.Cells.Item(2, 23).Value = 0.03217864529904591
$ws.Cells.Item(2, 24).Value = 0.03185898357847854
$ws.Cells.Item(2, 25).Value = 0.03168630631608825

# Row 3
$ws.Cells.Item(3, 3).Value = 0.8559999465942383
$ws.Cells.Item(3, 5).Value = 1715.578580270685
$ws.Cells.Item(3, 6).Value = 0.1763624076937847
$ws.Cells.Item(3, 7).Value = 0.1263645401470911
$ws.Cells.Item(3, 8).Value = 0.1020384455203502
$ws.Cells.Item(3, 9).Value = 0.08523151628714064
$ws.Cells.Item(3, 10).Value = 0.0658145945351216
$ws.Cells.Item(3, 11).Value = 0.06045834717086372
$ws.Cells.Item(3, 12).Value = 0.05372635611883469
$ws.Cells.Item(3, 13).Value = 0.04677246734473751
$ws.Cells.Item(3, 14).Value = 0.04609082253864018
$ws.Cells.Item(3, 15).Value = 0.04380612103035653
$ws.Cells.Item(3, 16).Value = 0.04046871622431129
$ws.Cells.Item(3, 17).Value = 0.03891525647602503
$ws.Cells.Item(3, 18).Value = 0.03747532082500028
$ws.Cells.Item(3, 19).Value = 0.03697736270137079
$ws.Cells.Item(3, 20).Value = 0.03583129126668937
$ws.Cells.Item(3, 21).Value = 0.0352052088692699
$ws.Cells.Item(3, 22).Value = 0.03443804797561788
$ws.Cells.Item(3, 23).Value = 0.03407361239690741
$ws.Cells.Item(3, 24).Value = 0.03374706320252515
$ws.Cells.Item(3, 25).Value = 0.03344207758812251

# Row 4
$ws.Cells.Item(4, 3).Value = 0.7579991817474365
$ws.Cells.Item(4, 5).Value = 1639.208327440214
$ws.Cells.Item(4, 6).Value = 0.1763624076937847
$ws.Cells.Item(4, 7).Value = 0.1348792417716657
$ws.Cells.Item(4, 8).Value = 0.1064780542695052
$ws.Cells.Item(4, 9).Value = 0.08338895958232009
$ws.Cells.Item(4, 10).Value = 0.07109310672687562
$ws.Cells.Item(4, 11).Value = 0.05707644621383016
$ws.Cells.Item(4, 12).Value = 0.05269301923888845
$ws.Cells.Item(4, 13).Value = 0.04720642890213662
$ws.Cells.Item(4, 14).Value = 0.04444040841295815
$ws.Cells.Item(4, 15).Value = 0.04250504918947473
$ws.Cells.Item(4, 16).Value = 0.03985298109276702
$ws.Cells.Item(4, 17).Value = 0.03826936320209191
$ws.Cells.Item(4, 18).Value = 0.03687013207744008
$ws.Cells.Item(4, 19).Value = 0.03539878650678455
$ws.Cells.Item(4, 20).Value = 0.03402232436309874
$ws.Cells.Item(4, 21).Value = 0.03366000579924027
$ws.Cells.Item(4, 22).Value = 0.03287234247776389
$ws.Cells.Item(4, 23).Value = 0.03207192759429546
$ws.Cells.Item(4, 24).Value = 0.03207192759429546
$ws.Cells.Item(4, 25).Value = 0.03195337870253828

# Row 5
$ws.Cells.Item(5, 3).Value = 0.8610007762908936
$ws.Cells.Item(5, 5).Value = 1598.631077141847
$ws.Cells.Item(5, 6).Value = 0.1763624076937847
$ws.Cells.Item(5, 7).Value = 0.1328641632549119
$ws.Cells.Item(5, 8).Value = 0.1047403048560694
$ws.Cells.Item(5, 9).Value = 0.08199541310190264
$ws.Cells.Item(5, 10).Value = 0.06655483234600262
$ws.Cells.Item(5, 11).Value = 0.05902362731007505
$ws.Cells.Item(5, 12).Value = 0.04842779973960514
$ws.Cells.Item(5, 13).Value = 0.04457124768618515
$ws.Cells.Item(5, 14).Value = 0.04250376263733432
$ws.Cells.Item(5, 15).Value = 0.03796772597101641
$ws.Cells.Item(5, 16).Value = 0.03681801582891269
$ws.Cells.Item(5, 17).Value = 0.03552024514080229
$ws.Cells.Item(5, 18).Value = 0.03510803604690538
$ws.Cells.Item(5, 19).Value = 0.03388264422017605
$ws.Cells.Item(5, 20).Value = 0.03321760881671676
$ws.Cells.Item(5, 21).Value = 0.03262180575220838
$ws.Cells.Item(5, 22).Value = 0.03230230445146789
$ws.Cells.Item(5, 23).Value = 0.03166300035594496
$ws.Cells.Item(5, 24).Value = 0.03130634213613796
$ws.Cells.Item(5, 25).Value = 0.03116239916455841

# Row 6
$ws.Cells.Item(6, 3).Value = 0.7569980621337891
$ws.Cells.Item(6, 5).Value = 1636.69141818678
$ws.Cells.Item(6, 6).Value = 0.1763624076937847
$ws.Cells.Item(6, 7).Value = 0.1257485336706867
$ws.Cells.Item(6, 8).Value = 0.09869983807965346
$ws.Cells.Item(6, 9).Value = 0.08087925441396365
$ws.Cells.Item(6, 10).Value = 0.06287755074402959
$ws.Cells.Item(6, 11).Value = 0.05750649716524688
$ws.Cells.Item(6, 12).Value = 0.04953384411395308
$ws.Cells.Item(6, 13).Value = 0.04420701251477521
$ws.Cells.Item(6, 14).Value = 0.04089653590020603
$ws.Cells.Item(6, 15).Value = 0.03937309346346302
$ws.Cells.Item(6, 16).Value = 0.03797555221528259
$ws.Cells.Item(6, 17).Value = 0.0370296314630482
$ws.Cells.Item(6, 18).Value = 0.0351324088088059
$ws.Cells.Item(6, 19).Value = 0.03448273109355746
$ws.Cells.Item(6, 20).Value = 0.03323039841595854
$ws.Cells.Item(6, 21).Value = 0.03285551410371019
$ws.Cells.Item(6, 22).Value = 0.03227327054130915
$ws.Cells.Item(6, 23).Value = 0.03218939289322043
$ws.Cells.Item(6, 24).Value = 0.03205666656184163
$ws.Cells.Item(6, 25).Value = 0.0319043161439918

# Row 7
$ws.Cells.Item(7, 3).Value = 0.7559986114501953
$ws.Cells.Item(7, 5).Value = 1652.738292180122
$ws.Cells.Item(7, 6).Value = 0.1763624076937847
$ws.Cells.Item(7, 7).Value = 0.1182437050952666
$ws.Cells.Item(7, 8).Value = 0.09548340358116639
$ws.Cells.Item(7, 9).Value = 0.08012567351658592
$ws.Cells.Item(7, 10).Value = 0.06639731938028642
$ws.Cells.Item(7, 11).Value = 0.05730215607895876
$ws.Cells.Item(7, 12).Value = 0.05449799746630783
$ws.Cells.Item(7, 13).Value = 0.05045240610184472
$ws.Cells.Item(7, 14).Value = 0.04412141374172823
$ws.Cells.Item(7, 15).Value = 0.04069530534904318
$ws.Cells.Item(7, 16).Value = 0.03832835765636929
$ws.Cells.Item(7, 17).Value = 0.03631017235813512
$ws.Cells.Item(7, 18).Value = 0.03502492252448234
$ws.Cells.Item(7, 19).Value = 0.03432939064918449
$ws.Cells.Item(7, 20).Value = 0.03393183971091718
$ws.Cells.Item(7, 21).Value = 0.03335786210659412
$ws.Cells.Item(7, 22).Value = 0.03248014786489616
$ws.Cells.Item(7, 23).Value = 0.03248014786489616
$ws.Cells.Item(7, 24).Value = 0.03248014786489616
$ws.Cells.Item(7, 25).Value = 0.03221712070526552

# Row 8
$ws.Cells.Item(8, 3).Value = 0.860037088394165
$ws.Cells.Item(8, 5).Value = 1509.654657218398
$ws.Cells.Item(8, 6).Value = 0.1763624076937847
$ws.Cells.Item(8, 7).Value = 0.1202340009318255
$ws.Cells.Item(8, 8).Value = 0.09833797571686118
$ws.Cells.Item(8, 9).Value = 0.07315512528984502
$ws.Cells.Item(8, 10).Value = 0.06189495683590591
$ws.Cells.Item(8, 11).Value = 0.05276503530350155
$ws.Cells.Item(8, 12).Value = 0.04596035170285108
$ws.Cells.Item(8, 13).Value = 0.03949592627949728
$ws.Cells.Item(8, 14).Value = 0.03784797821382108
$ws.Cells.Item(8, 15).Value = 0.0351078621414107
$ws.Cells.Item(8, 16).Value = 0.03393091595358037
$ws.Cells.Item(8, 17).Value = 0.03314490568814368
$ws.Cells.Item(8, 18).Value = 0.03258828138589678
$ws.Cells.Item(8, 19).Value = 0.03132083658371653
$ws.Cells.Item(8, 20).Value = 0.03120214404909025
$ws.Cells.Item(8, 21).Value = 0.03062648298374655
$ws.Cells.Item(8, 22).Value = 0.03018361314833
$ws.Cells.Item(8, 23).Value = 0.02974871747945112
$ws.Cells.Item(8, 24).Value = 0.02963673536583354
$ws.Cells.Item(8, 25).Value = 0.02942796602764908

# Row 9
$ws.Cells.Item(9, 3).Value = 0.7769944667816162
$ws.Cells.Item(9, 5).Value = 1629.17913074275
$ws.Cells.Item(9, 6).Value = 0.1763624076937847
$ws.Cells.Item(9, 7).Value = 0.1241050196672212
$ws.Cells.Item(9, 8).Value = 0.09468685151270141
$ws.Cells.Item(9, 9).Value = 0.07846231750126748
$ws.Cells.Item(9, 10).Value = 0.06215609006798381
$ws.Cells.Item(9, 11).Value = 0.05555157635751428
$ws.Cells.Item(9, 12).Value = 0.05027196748475667
$ws.Cells.Item(9, 13).Value = 0.0430860783254331
$ws.Cells.Item(9, 14).Value = 0.04219049019967016
$ws.Cells.Item(9, 15).Value = 0.03973941329316557
$ws.Cells.Item(9, 16).Value = 0.03772709994819048
$ws.Cells.Item(9, 17).Value = 0.03629713363698182
$ws.Cells.Item(9, 18).Value = 0.03515779587272554
$ws.Cells.Item(9, 19).Value = 0.03386573739978496
$ws.Cells.Item(9, 20).Value = 0.03342168569798855
$ws.Cells.Item(9, 21).Value = 0.0333357501383874
$ws.Cells.Item(9, 22).Value = 0.03271707198430868
$ws.Cells.Item(9, 23).Value = 0.03225136547663787
$ws.Cells.Item(9, 24).Value = 0.03201794039697078
$ws.Cells.Item(9, 25).Value = 0.03175787779225633

# Row 10
$ws.Cells.Item(10, 3).Value = 0.7090342044830322
$ws.Cells.Item(10, 5).Value = 1738.884846641593
$ws.Cells.Item(10, 6).Value = 0.1763624076937847
$ws.Cells.Item(10, 7).Value = 0.1341004092578325
$ws.Cells.Item(10, 8).Value = 0.09808143689876682
$ws.Cells.Item(10, 9).Value = 0.08121672649955827
$ws.Cells.Item(10, 10).Value = 0.06333204452698174
$ws.Cells.Item(10, 11).Value = 0.05842157055022604
$ws.Cells.Item(10, 12).Value = 0.05211522471589179
$ws.Cells.Item(10, 13).Value = 0.04565099832168074
$ws.Cells.Item(10, 14).Value = 0.04425244501595477
$ws.Cells.Item(10, 15).Value = 0.04077093980750659
$ws.Cells.Item(10, 16).Value = 0.04008316138651232
$ws.Cells.Item(10, 17).Value = 0.03935589456409404
$ws.Cells.Item(10, 18).Value = 0.03888512717362943
$ws.Cells.Item(10, 19).Value = 0.03780155422777722
$ws.Cells.Item(10, 20).Value = 0.03669165015705624
$ws.Cells.Item(10, 21).Value = 0.03537375049584821
$ws.Cells.Item(10, 22).Value = 0.03531904763530332
$ws.Cells.Item(10, 23).Value = 0.03464552953318453
$ws.Cells.Item(10, 24).Value = 0.03421141518570118
$ws.Cells.Item(10, 25).Value = 0.03389639077274059

# Row 11
$ws.Cells.Item(11, 3).Value = 0.7389984130859375
$ws.Cells.Item(11, 5).Value = 1619.547529769727
$ws.Cells.Item(11, 6).Value = 0.1763624076937847
$ws.Cells.Item(11, 7).Value = 0.1192790818141171
$ws.Cells.Item(11, 8).Value = 0.09458245683278062
$ws.Cells.Item(11, 9).Value = 0.08004368783251962
$ws.Cells.Item(11, 10).Value = 0.07279067849665462
$ws.Cells.Item(11, 11).Value = 0.06306303747011033
$ws.Cells.Item(11, 12).Value = 0.05516711879972319
$ws.Cells.Item(11, 13).Value = 0.04993685499942606
$ws.Cells.Item(11, 14).Value = 0.04603354970106819
$ws.Cells.Item(11, 15).Value = 0.04060754762103787
$ws.Cells.Item(11, 16).Value = 0.03880228962400269
$ws.Cells.Item(11, 17).Value = 0.03666663948976245
$ws.Cells.Item(11, 18).Value = 0.03485637882208108
$ws.Cells.Item(11, 19).Value = 0.03470775699119931
$ws.Cells.Item(11, 20).Value = 0.03348187890611519
$ws.Cells.Item(11, 21).Value = 0.03281139778005213
$ws.Cells.Item(11, 22).Value = 0.03204032492713067
$ws.Cells.Item(11, 23).Value = 0.03203920798490185
$ws.Cells.Item(11, 24).Value = 0.03169549392742469
$ws.Cells.Item(11, 25).Value = 0.03157012728595958

